# The domain.xlsx import template previously generated two extra header
# columns for the "is_locked" and "is_enabled" boolean fields (with their
# accompanying data-validation generator snippets). Those two columns
# (C:D -> is_locked_lbl, is_enabled_lbl) are removed, and the remaining
# columns (order_by, rem) shift left to take their place.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1:D1").EntireColumn.Delete()
